$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for D/E columns so numeric-looking strings are not
# auto-converted to numbers by Excel (matches original inlineStr text cells).
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '43.725.92'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '2.286.25'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '115.52'
$ws.Range('E5').Value = '  +12.77%  '
$ws.Range('D6').Value = '268.41'
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('D7').Value = '0.627'
$ws.Range('E7').Value = '  +0.83%  '
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('D10').Value = '48.89'
$ws.Range('E10').Value = '  +8.08%  '
$ws.Range('D11').Value = '0.0943'
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('D12').Value = '8.88'
$ws.Range('E12').Value = '  +11.46%  '
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('D14').Value = '15.81'
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = '2.638.13'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '0.881'
$ws.Range('E16').Value = '  +2.93%  '
$ws.Range('D17').Value = '2.281.43'
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('D18').Value = '43.686.18'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('E19').Value = '  -1.64%  '
$ws.Range('D20').Value = '6.99'
$ws.Range('E20').Value = '  +11.86%  '
$ws.Range('D21').Value = '72.21'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').Value = '2.41'
$ws.Range('E22').Value = '  -3.15%  '
$ws.Range('D23').Value = '9.86'
$ws.Range('E23').Value = '  +7.76%  '
$ws.Range('D24').Value = '232.63'
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('D25').Value = '2.92'
$ws.Range('E25').Value = '  +1.90%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '11.63'
$ws.Range('E27').Value = '  +3.79%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').Value = '42.15'
$ws.Range('E28').Value = '  +6.08%  '
$ws.Range('B29').Value = 'WEMIXToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D29').Value = '3.39'
$ws.Range('E29').Value = '  -1.78%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '2.24'
$ws.Range('E30').Value = '  -2.35%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '173.09'
$ws.Range('E31').Value = '  -2.42%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '21.58'
$ws.Range('E32').Value = '  -1.04%  '
$ws.Range('D33').Value = '0.0928'
$ws.Range('E33').Value = '  +3.21%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '5.70'
$ws.Range('E34').Value = '  +4.66%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').Value = '0.127'
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '4.67'
$ws.Range('E36').Value = '  -4.01%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.0357'
$ws.Range('E37').Value = '  +1.16%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '0.107'
$ws.Range('E38').Value = '  -1.41%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = '3.79'
$ws.Range('E39').Value = '  +6.05%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '14.52'
$ws.Range('E40').Value = '  +18.46%  '
$ws.Range('B41').Value = 'MultiversX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D41').Value = '74.42'
$ws.Range('E41').Value = '  +14.12%  '
$ws.Range('E42').Value = '  +3.86%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '0.241'
$ws.Range('E43').Value = '  +2.26%  '
$ws.Range('B44').Value = 'THORChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D44').Value = '6.36'
$ws.Range('E44').Value = '  +21.28%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').Value = '1.39'
$ws.Range('E46').Value = '  -0.47%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '8.70'
$ws.Range('E47').Value = '  -0.87%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').Value = '1.26'
$ws.Range('E48').Value = '  +4.15%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '102.79'
$ws.Range('E49').Value = '  +4.38%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.100'
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '0.458'
$ws.Range('E51').Value = '  +1.90%  '

# Restore default (unstyled) cell style now that values are written as text,
# so the saved cells don't carry a stray style index like the source file.
$ws.Range('D2:E51').Style = 'Normal'

